# Swap the display order of slide 9 ("Authorization") and slide 10
# ("Handling Authorization With Hashing") in the deck.
#
# Before: position 9  = "Authorization" (simple content slide)
#         position 10 = "Handling Authorization With Hashing" (designed slide)
# After:  position 9  = "Handling Authorization With Hashing"
#         position 10 = "Authorization"
#
# Moving the slide that is currently at index 10 to index 9 pushes the
# former slide 9 down to index 10, producing the desired reordering.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$s.MoveTo(9)
